# [AS]- Done Refactoring in MCA, UserDetailsPS, SearchOnlineAccount_2 and LifeSpanFlow
#
# Adds a new "Install - POST-PAY" test-case block (group header + 9 rows) to
# the end of the MCA.CreateUpdateMsisdnProfile sheet, mirroring the existing
# "Install - PRE-PAY" block (rows 77:85) content, and updates the saved
# window/view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MCA.CreateUpdateMsisdnProfile")

# --- 1. Stamp the formatting for the new block by copying existing,
#        identically-styled rows down onto the new row range. -----------

# Group-header row (style pattern used by every section banner, e.g. A116:I116)
$ws.Range("A128:I128").Merge() | Out-Null
$ws.Range("A116:I116").Copy() | Out-Null
$ws.Range("A128:I128").PasteSpecial(-4122) | Out-Null

# Nine data rows - reuse the plain "105pt" row formatting (A118:I126)
$ws.Range("A118:I126").Copy() | Out-Null
$ws.Range("A129:I137").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- 2. Row heights: first data row of the block is taller (wraps more
#        text), the rest share the standard height. ----------------------

$ws.Rows.Item(128).RowHeight = 15.75
$ws.Rows.Item(129).RowHeight = 165
$ws.Rows.Item(130).RowHeight = 105
$ws.Rows.Item(131).RowHeight = 105
$ws.Rows.Item(132).RowHeight = 105
$ws.Rows.Item(133).RowHeight = 105
$ws.Rows.Item(134).RowHeight = 105
$ws.Rows.Item(135).RowHeight = 105
$ws.Rows.Item(136).RowHeight = 105
$ws.Rows.Item(137).RowHeight = 105

# --- 3. Values. Re-use the text already stored in the mirror "Install -
#        PRE-PAY" block (rows 77:85) so the strings are transcribed
#        exactly, rather than re-typed. ----------------------------------

$ws.Range("A128").Value = "Install - POST-PAY"

$srcRows = 77, 78, 79, 80, 81, 82, 83, 84, 85
$dstRows = 129, 130, 131, 132, 133, 134, 135, 136, 137
$ids = 112, 113, 114, 115, 116, 117, 118, 119, 120

for ($i = 0; $i -lt $srcRows.Length; $i++) {
    $srcRow = $srcRows[$i]
    $dstRow = $dstRows[$i]

    $ws.Cells.Item($dstRow, 1).Value = $ids[$i]
    $ws.Cells.Item($dstRow, 2).Value = $ws.Cells.Item($srcRow, 2).Value2
    $ws.Cells.Item($dstRow, 3).Value = $ws.Cells.Item($srcRow, 3).Value2
    $ws.Cells.Item($dstRow, 9).Value = "Y"
}

# --- 4. View state: scroll/selection as left by the edit ----------------

$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$ws.Range("C130").Select()

$excel.ActiveWindow.WindowState = $excel.ActiveWindow.WindowState
$excel.Width = 15345
$excel.Height = 4560
